$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", copying the header style (s="1") from E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for data rows 2-67 (no special style, matches other data cells)
$timestamps = @(
    "2021-10-05 13:42:24.580779",
    "2021-10-05 13:42:24.580791",
    "2021-10-05 13:42:24.580795",
    "2021-10-05 13:42:24.580798",
    "2021-10-05 13:42:24.580801",
    "2021-10-05 13:42:24.580805",
    "2021-10-05 13:42:24.580807",
    "2021-10-05 13:42:24.580810",
    "2021-10-05 13:42:24.580814",
    "2021-10-05 13:42:24.580817",
    "2021-10-05 13:42:24.580820",
    "2021-10-05 13:42:24.580823",
    "2021-10-05 13:42:24.580826",
    "2021-10-05 13:42:24.580829",
    "2021-10-05 13:42:24.580832",
    "2021-10-05 13:42:24.580834",
    "2021-10-05 13:42:24.580838",
    "2021-10-05 13:42:24.580841",
    "2021-10-05 13:42:24.580844",
    "2021-10-05 13:42:24.580847",
    "2021-10-05 13:42:24.580850",
    "2021-10-05 13:42:24.580853",
    "2021-10-05 13:42:24.580855",
    "2021-10-05 13:42:24.580858",
    "2021-10-05 13:42:24.580862",
    "2021-10-05 13:42:24.580865",
    "2021-10-05 13:42:24.580868",
    "2021-10-05 13:42:24.580870",
    "2021-10-05 13:42:24.580873",
    "2021-10-05 13:42:24.580876",
    "2021-10-05 13:42:24.580879",
    "2021-10-05 13:42:24.580882",
    "2021-10-05 13:42:24.580886",
    "2021-10-05 13:42:24.580889",
    "2021-10-05 13:42:24.580892",
    "2021-10-05 13:42:24.580894",
    "2021-10-05 13:42:24.580897",
    "2021-10-05 13:42:24.580900",
    "2021-10-05 13:42:24.580903",
    "2021-10-05 13:42:24.580906",
    "2021-10-05 13:42:24.580909",
    "2021-10-05 13:42:24.580912",
    "2021-10-05 13:42:24.580915",
    "2021-10-05 13:42:24.580918",
    "2021-10-05 13:42:24.580921",
    "2021-10-05 13:42:24.580924",
    "2021-10-05 13:42:24.580927",
    "2021-10-05 13:42:24.580930",
    "2021-10-05 13:42:24.580933",
    "2021-10-05 13:42:24.580936",
    "2021-10-05 13:42:24.580939",
    "2021-10-05 13:42:24.580942",
    "2021-10-05 13:42:24.580945",
    "2021-10-05 13:42:24.580948",
    "2021-10-05 13:42:24.580951",
    "2021-10-05 13:42:24.580954",
    "2021-10-05 13:42:24.580957",
    "2021-10-05 13:42:24.580960",
    "2021-10-05 13:42:24.580963",
    "2021-10-05 13:42:24.580966",
    "2021-10-05 13:42:24.580969",
    "2021-10-05 13:42:24.580972",
    "2021-10-05 13:42:24.580975",
    "2021-10-05 13:42:24.580978",
    "2021-10-05 13:42:24.580982",
    "2021-10-05 13:42:24.580985"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

